# New non-convex experiment values: shift x down by 1.65 and y up by 1.65,
# then update all the dependent expression-evaluation strings accordingly.
# All of these cells hold text (not numeric) values in the workbook, so we
# force a text number-format before writing, and reset the style back to
# Normal afterwards so no stray cell-style index gets introduced.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$value)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# NOTE: worksheet names are matched case-insensitively by Worksheets.Item(name),
# and this workbook has both "Vector_bf" and "Vector_BF" sheets, so we must
# address sheets by their (1-based) position instead of by name.

# --- Restricciones_del_lider (sheet #2) ---
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2 "A2" "2.8499999999999996 - x"
Set-TextValue $ws2 "B2" "-3.3499999999999996"
Set-TextValue $ws2 "D2" "0.3"
Set-TextValue $ws2 "A3" "-2.8499999999999996 + x"
Set-TextValue $ws2 "B3" "2.3499999999999996"
Set-TextValue $ws2 "D3" "0.09"

# --- Restricciones_del_follower (sheet #3) ---
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3 "A2" "4.449999999999999 - y"
Set-TextValue $ws3 "B2" "-5.449999999999999"
Set-TextValue $ws3 "D2" "0.19"
Set-TextValue $ws3 "E2" "0"
Set-TextValue $ws3 "F2" "2.1"
Set-TextValue $ws3 "A3" "-4.449999999999999 + y"
Set-TextValue $ws3 "B3" "3.4499999999999993"
Set-TextValue $ws3 "D3" "0.79"
Set-TextValue $ws3 "E3" "0"
Set-TextValue $ws3 "F3" "9.3"

# --- Punto_modificado (sheet #4) ---
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4 "A2" "2.8499999999999996"
Set-TextValue $ws4 "B2" "4.449999999999999"

# --- Vector_bf (sheet #5) ---
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5 "A2" "-3.4499999999999997"

# --- Vector_BF (sheet #6) ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6 "A2" "1.21"
Set-TextValue $ws6 "A3" "1.0"
